$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.276.10"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.307.80"
$ws.Range("E3").Value = "  +4.39%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.03"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.643"
$ws.Range("E6").Value = "  +3.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.54"
$ws.Range("E7").Value = "  +10.14%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +6.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.54"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0994"
$ws.Range("E11").Value = "  +6.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.40"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.52"
$ws.Range("E13").Value = "  +7.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.653.30"
$ws.Range("E15").Value = "  +4.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.62"
$ws.Range("E16").Value = "  +7.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.885"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.311.69"
$ws.Range("E18").Value = "  +4.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.120.78"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("E20").Value = "  +5.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.35"
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.89"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.32"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("E24").Value = "  +10.36%  "
$ws.Range("E26").Value = "  +4.88%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.62"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.20"
$ws.Range("E32").Value = "  +3.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.42"
$ws.Range("E33").Value = "  +8.55%  "
$ws.Range("E34").Value = "  +6.02%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.20"
$ws.Range("E35").Value = "  +23.66%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0816"
$ws.Range("E36").Value = "  +5.11%  "
$ws.Range("E37").Value = "  +4.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.60"
$ws.Range("E38").Value = "  +12.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.80"
$ws.Range("E39").Value = "  +4.70%  "
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  +21.42%  "
$ws.Range("E42").Value = "  +6.49%  "
$ws.Range("E43").Value = "  +6.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.216"
$ws.Range("E44").Value = "  +10.85%  "
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.45"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.18"
$ws.Range("E46").Value = "  +7.55%  "
$ws.Range("E47").Value = "  -5.17%  "
$ws.Range("E48").Value = "  +4.55%  "
$ws.Range("E49").Value = "  +3.66%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "99.01"
$ws.Range("E51").Value = "  +6.89%  "
